$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function DeleteHyperlinkAt($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        $a = $h.Range.Address()
        if ($a -eq $addr) {
            $h.Delete()
            return
        }
    }
}

# Row 3 (com.hamxa.shaynachim / bitcoin guide / eligitel@gmail.com / ronenchen27@gmail.com / ...)
# was removed from the review table; every row below it shifts up one.
$ws.Rows("3:3").Delete()

# Row deletion does not retarget the worksheet's <hyperlinks>, so the stale
# entries (old C3/D3, which pointed at the now-deleted eligitel/ronenchen
# addresses) and the old C4 entry (armonravid, now sitting in row 3) need to
# be fixed up by hand.
DeleteHyperlinkAt $ws '$D$3'
DeleteHyperlinkAt $ws '$C$3'
DeleteHyperlinkAt $ws '$C$4'

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:armonravid@gmail.com", "", "", "armonravid@gmail.com")

# Adding a hyperlink pushes in Excel's built-in "Hyperlink" cell style; put
# the cell's formatting back the way it was (same font/alignment as the
# other email cells in the table).
$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").Font.Size = 11
$ws.Range("C3").Font.Color = 0
$ws.Range("C3").Font.Underline = 0
$ws.Range("C3").HorizontalAlignment = -4108

# Adding the hyperlink also registers a new built-in "Hyperlink" named cell
# style on the workbook; since C3 no longer uses it, drop it again so the
# style table doesn't grow.
$wb.Styles.Item("Hyperlink").Delete()

# Match the saved selection state (A3) left behind by the edit.
$null = $ws.Range("A3").Select()
